{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// This applies the \"Update symptom checker docs for 1/20 changes\" edit:\n//  1. Remove the \"IMAGE\" placeholder block (the bold \"IMAGE\" label, the\n//     \"Anatomy of the eye IM02853\" title line, its caption paragraph and\n//     the blank spacer paragraphs around them) that used to sit right\n//     before the \"More Information\" section, collapsing it so the\n//     paragraph that used to read \"IMAGE\" now reads \"More Information\"\n//     (re-using its existing bold-run formatting), and the paragraph\n//     that used to hold the original \"More Information\" text is removed.\n//  2. Bump the trailing page marker from \" PAGE 4\" to \" PAGE 5\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the paragraph that says exactly \"IMAGE\" and, after it, the next\n// paragraph that says exactly \"More Information\".\nlet imageIndex = -1;\nlet moreInfoIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (imageIndex === -1 && items[i].text === \"IMAGE\") {\n    imageIndex = i;\n    continue;\n  }\n  if (imageIndex !== -1 && moreInfoIndex === -1 && items[i].text === \"More Information\") {\n    moreInfoIndex = i;\n    break;\n  }\n}\n\nif (imageIndex !== -1 && moreInfoIndex !== -1) {\n  // Delete every paragraph from the original \"More Information\" paragraph\n  // back through (but not including) the \"IMAGE\" paragraph. Walking\n  // backwards keeps the remaining indices valid as we go.\n  for (let i = moreInfoIndex; i > imageIndex; i--) {\n    items[i].delete();\n  }\n  await context.sync();\n\n  // Re-purpose the run that used to say \"IMAGE\" so it says\n  // \"More Information\" instead (its bold formatting carries over).\n  items[imageIndex].getRange().insertText(\"More Information\", \"Replace\");\n  await context.sync();\n}\n\n// \" PAGE 4\" -> \" PAGE 5\"\nconst pageResults = body.search(\" PAGE 4\", { matchCase: true });\nawait context.sync();\nif (pageResults.items.length > 0) {\n  pageResults.items[0].insertText(\" PAGE 5\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d below.\n#\n# This applies the \"Update symptom checker docs for 1/20 changes\" edit:\n#  1. Remove the \"IMAGE\" placeholder block (the bold \"IMAGE\" label, the\n#     \"Anatomy of the eye IM02853\" title line, its caption paragraph and\n#     the blank spacer paragraphs around them) that used to sit right\n#     before the \"More Information\" section, collapsing it so the\n#     paragraph that used to read \"IMAGE\" now reads \"More Information\"\n#     (re-using its existing bold-run formatting), and the paragraph\n#     that used to hold the original \"More Information\" text is removed.\n#  2. Bump the trailing page marker from \" PAGE 4\" to \" PAGE 5\".\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n# Locate the paragraph that reads exactly \"IMAGE\" and, after it, the next\n# paragraph that reads exactly \"More Information\".\n$imageIndex = -1\n$moreInfoIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $t = $paras.Item($i).Range.Text.TrimEnd([char]13)\n    if ($imageIndex -eq -1) {\n        if ($t -eq \"IMAGE\") {\n            $imageIndex = $i\n        }\n        continue\n    }\n    if ($t -eq \"More Information\") {\n        $moreInfoIndex = $i\n        break\n    }\n}\n\nif ($imageIndex -ne -1 -and $moreInfoIndex -ne -1) {\n    # Delete every paragraph from the original \"More Information\" paragraph\n    # back through (but not including) the \"IMAGE\" paragraph. Deleting from\n    # the highest index down keeps the lower indices valid.\n    for ($i = $moreInfoIndex; $i -gt $imageIndex; $i--) {\n        $paras.Item($i).Range.Delete()\n    }\n\n    # Re-purpose the run that used to say \"IMAGE\" so it says\n    # \"More Information\" instead (its bold formatting carries over).\n    $imgRange = $paras.Item($imageIndex).Range\n    $textRange = $d.Range($imgRange.Start, $imgRange.Start + 5)\n    $textRange.Text = \"More Information\"\n}\n\n# \" PAGE 4\" -> \" PAGE 5\"\n$find = $d.Content.Find\n$find.Text = \" PAGE 4\"\n$find.Replacement.Text = \" PAGE 5\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n"}
